$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column J mirrors column I (LO phase noise @ values) for rows 1-3
$ws.Range("J1").Value = $ws.Range("I1").Value2
$ws.Range("J2").Value = $ws.Range("I2").Value2
$ws.Range("J3").Value = $ws.Range("I3").Value2

# Apply an explicit "General" number format to G7 (new style entry)
$ws.Range("G7").NumberFormat = "General"

# Update the active selection/view to H11
$ws.Range("H11").Select() | Out-Null
